$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 38 (pushes "Ken" and everything below it down by one),
# matching the author inserting "Kazuya" into the character attribute table.
$ws.Rows.Item(38).Insert()

# Populate the newly inserted row with Kazuya's attribute data.
$ws.Cells.Item(38, 1).Value = "Kazuya"
$ws.Cells.Item(38, 2).Value = 1
$ws.Cells.Item(38, 3).Value = 1
$ws.Cells.Item(38, 4).Value = 133
$ws.Cells.Item(38, 5).Value = 0.108
$ws.Cells.Item(38, 6).Value = 2.72
$ws.Cells.Item(38, 7).Value = 0.079968
$ws.Cells.Item(38, 8).Value = 1.8
$ws.Cells.Item(38, 9).Value = 0.115

# Mirror the author's final cursor position (sheetView selection) from the diff.
$ws.Range("K33").Select()
